$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Codelist Usage")

# Remove the duplicate "ThermalPipe / thermalProductType" rows.
# Row 130 is a duplicate of row 129 (ThermalProductTypeIMKLValue) - delete first (higher row).
$ws.Rows.Item(130).Delete()
# Row 128 is a duplicate of row 127 (ThermalProductTypeExtendedValue).
$ws.Rows.Item(128).Delete()
